$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting of the row immediately above (row 10, B:J) down onto
# the new row 11 first, so the new row's cells pick up the same styles
# (date format, borders, wrap text, etc.) that rows 9/10 already use.
$ws.Range("B10:J10").Copy()
$ws.Range("B11:J11").PasteSpecial(-4122)
$ws.Rows("11:11").RowHeight = 38.25

# New risk item (row 11 of the risk list table)
$ws.Range("B11").Value = 41739
$ws.Range("C11").Value = "Risco da falta de domínio nos padrões de UML pode gerar problemas na equipe no desenvolvimento da aplicação"
$ws.Range("D11").Value = "Realizar treinamento sobre o padrão UMLcom todos os integrantes da equipe de desenvolvimento "
$ws.Range("E11").Value = "D"
$ws.Range("F11").Value = 3
$ws.Range("G11").Value = 0.3
$ws.Range("H11").Formula = "=+F11*G11"
$ws.Range("I11").Value = "Analista de Desenvolvimento"
$ws.Range("J11").Value = "Organizar os treinamentos aos Sábados para não gerar aumento de custo significativo ao projeto"

# Match the author's final selection (cursor moved one row past the new data)
$ws.Range("J12").Select() | Out-Null
